# Update api endpoint validation
# Appends one new telemetry row to each of the four worksheets, matching
# the new rows captured by the validation run.

$wb = $excel.ActiveWorkbook

# --- Sheet "ROW50-FE-LIFTER" (sheet 1): new row 26 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 26
$ws1.Cells.Item($r, 1).Value = 45737.13607120371
$ws1.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x7a"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws1.Cells.Item($r, 8).Value = 378
$ws1.Cells.Item($r, 9).Value = 14

# --- Sheet "ROW50-MID-LIFTER" (sheet 2): new row 28 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 28
$ws2.Cells.Item($r, 1).Value = 45737.11314814815
$ws2.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x82"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400
# Keep this long digit string as text (matches source data, avoids float rounding)
$ws2.Cells.Item($r, 7).Value = "'568631262647113771663628"
$ws2.Cells.Item($r, 8).Value = 386
$ws2.Cells.Item($r, 9).Value = 25

# --- Sheet "ROW11-FE-LIFTER" (sheet 3): new row 26 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 26
$ws3.Cells.Item($r, 1).Value = 45737.16007710648
$ws3.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x7a"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws3.Cells.Item($r, 8).Value = 378
$ws3.Cells.Item($r, 9).Value = 20

# --- Sheet "ROW11-MID-LIFTER" (sheet 4): new row 26 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 26
$ws4.Cells.Item($r, 1).Value = 45737.30617762732
$ws4.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x82"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws4.Cells.Item($r, 8).Value = 386
$ws4.Cells.Item($r, 9).Value = 25

Write-Host "Appended rows to all four sheets"
